$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-07 Friday" "2025-11-08 Saturday"

Replace-Text "83×47=3901" "93×95=8835"
Replace-Text "79×38=3002" "36×35=1260"
Replace-Text "70×77=5390" "62×46=2852"
Replace-Text "76×74=5624" "91×33=3003"
Replace-Text "11×49=539" "88×55=4840"
Replace-Text "23×42=966" "47×32=1504"
Replace-Text "65×19=1235" "88×85=7480"
Replace-Text "27×58=1566" "42×93=3906"
Replace-Text "75×56=4200" "40×96=3840"
Replace-Text "92×76=6992" "73×11=803"
Replace-Text "53×92=4876" "59×51=3009"
Replace-Text "75×35=2625" "20×83=1660"
Replace-Text "44×26=1144" "40×75=3000"
Replace-Text "35×30=1050" "24×82=1968"
Replace-Text "18×74=1332" "39×25=975"
Replace-Text "88×81=7128" "37×33=1221"
Replace-Text "98×79=7742" "84×71=5964"
Replace-Text "23×49=1127" "93×34=3162"
Replace-Text "96×41=3936" "40×21=840"
Replace-Text "67×22=1474" "63×97=6111"
Replace-Text "91×62=5642" "85×70=5950"
Replace-Text "16×53=848" "89×47=4183"
Replace-Text "62×54=3348" "17×43=731"
Replace-Text "30×40=1200" "55×73=4015"
Replace-Text "16×18=288" "15×83=1245"

Write-Output "Done"
